# Adds a dropdown (data validation) to switch between different days, and
# stores each day's task text in a dedicated column (A, C, E) on row 1 so a
# "save task for each day" helper can retrieve/restore it later.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the task text for "day 1" (column A) which used to be blank.
$ws.Range("A1").Value = "Read Chapters 21-23"

# New hidden/helper columns holding the task text for the other days that
# the dropdown can switch between.
$ws.Range("C1").Value = "Read Chapters 24-25"
$ws.Range("E1").Value = "Read Chapters 26-29"

# Column E mirrors column A's per-row blank placeholders (single space)
# for rows 2-5, matching the existing pattern used by columns A/B.
$ws.Range("E2").Value = " "
$ws.Range("E3").Value = " "
$ws.Range("E4").Value = " "
$ws.Range("E5").Value = " "

# Select C2 - the cell just below the new "day" dropdown header - as the
# active selection, matching where the user left off after adding it.
$ws.Range("C2").Select()
